$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.420.86"
$ws.Range("E2").Value = "  -1.56%  "
$ws.Range("D3").Value = "2.483.15"
$ws.Range("E3").Value = "  -1.84%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'525.08"
$ws.Range("E5").Value = "  -3.10%  "
$ws.Range("D6").Value = "'133.65"
$ws.Range("E6").Value = "  -3.65%  "
$ws.Range("D8").Value = "'0.560"
$ws.Range("E8").Value = "  -1.16%  "
$ws.Range("D9").Value = "'0.0999"
$ws.Range("E9").Value = "  -1.77%  "
$ws.Range("E10").Value = "  -1.88%  "
$ws.Range("D11").Value = "'5.43"
$ws.Range("E11").Value = "  +0.85%  "
$ws.Range("D12").Value = "'0.343"
$ws.Range("E12").Value = "  -2.05%  "
$ws.Range("D13").Value = "2.925.36"
$ws.Range("E13").Value = "  -1.71%  "
$ws.Range("D14").Value = "58.384.56"
$ws.Range("E14").Value = "  -1.51%  "
$ws.Range("D15").Value = "'22.42"
$ws.Range("E15").Value = "  -3.77%  "
$ws.Range("E16").Value = "  -2.26%  "
$ws.Range("D17").Value = "2.486.54"
$ws.Range("E17").Value = "  -1.69%  "
$ws.Range("D18").Value = "'10.93"
$ws.Range("E18").Value = "  -2.00%  "
$ws.Range("E19").Value = "  -2.46%  "
$ws.Range("D20").Value = "'321.55"
$ws.Range("E20").Value = "  -1.59%  "
$ws.Range("D21").Value = "'0.998"
$ws.Range("E21").Value = "  -0.18%  "
$ws.Range("D22").Value = "'5.81"
$ws.Range("E22").Value = "  -2.50%  "
$ws.Range("D23").Value = "'64.37"
$ws.Range("E23").Value = "  -1.58%  "
$ws.Range("E24").Value = "  -3.11%  "
$ws.Range("D25").Value = "'1.00"
$ws.Range("D26").Value = "'0.161"
$ws.Range("E26").Value = "  -3.46%  "
$ws.Range("E27").Value = "  -2.90%  "
$ws.Range("D28").Value = "0.0₃0752"
$ws.Range("E28").Value = "  -3.89%  "
$ws.Range("E29").Value = "  -4.94%  "
$ws.Range("E30").Value = "  -4.85%  "
$ws.Range("D31").Value = "'166.33"
$ws.Range("E31").Value = "  -1.29%  "
$ws.Range("E32").Value = "  -5.70%  "
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("D35").Value = "'18.24"
$ws.Range("E35").Value = "  -1.61%  "
$ws.Range("E36").Value = "  -8.97%  "
$ws.Range("E37").Value = "  -4.02%  "
$ws.Range("E38").Value = "  -4.32%  "
$ws.Range("D39").Value = "'0.798"
$ws.Range("E39").Value = "  -3.49%  "
$ws.Range("D40").Value = "'3.53"
$ws.Range("E40").Value = "  -3.67%  "
$ws.Range("D41").Value = "'278.31"
$ws.Range("E41").Value = "  -2.38%  "
$ws.Range("E42").Value = "  -5.70%  "
$ws.Range("D43").Value = "'0.595"
$ws.Range("E43").Value = "  -2.36%  "
$ws.Range("D44").Value = "'127.71"
$ws.Range("E44").Value = "  -2.89%  "
$ws.Range("D45").Value = "'0.0913"
$ws.Range("E45").Value = "  -2.40%  "
$ws.Range("D46").Value = "'0.0496"
$ws.Range("E46").Value = "  -3.02%  "
$ws.Range("E47").Value = "  -2.83%  "
$ws.Range("D48").Value = "'17.25"
$ws.Range("E48").Value = "  -1.67%  "
$ws.Range("D49").Value = "1.741.61"
$ws.Range("E49").Value = "  -1.35%  "
$ws.Range("E50").Value = "  -1.68%  "
$ws.Range("E51").Value = "  -2.12%  "
